$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 0.176022
$ws.Range("H2").Value = 0.5280659999999999
$ws.Range("I2").Value = 0.03293066697281707
$ws.Range("J2").Value = 0.03293066697281707
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 63.55492266666666
$ws.Range("N2").Value = 190.664768
$ws.Range("O2").Value = 0.9936031556622397
$ws.Range("P2").Value = 0.9936031556622397
$ws.Range("Q2").Value = 11.187064597632
$ws.Range("R2").Value = 100.683581378688
$ws.Range("S2").Value = 0.03272001462225334
$ws.Range("T2").Value = 0.03272001462225333
# Row 3
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 0.176022
$ws.Range("H3").Value = 0.5280659999999999
$ws.Range("I3").Value = 0.03293066697281707
$ws.Range("J3").Value = 0.03293066697281707
$ws.Range("O3").Value = 0.000612609346703606
$ws.Range("P3").Value = 0.000612609346703606
$ws.Range("Q3").Value = 0.00689742207
$ws.Range("R3").Value = 0.06207679862999999
$ws.Range("S3").Value = 0.00002017363438073148
$ws.Range("T3").Value = 0.00002017363438073148
# Row 4
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 0.176022
$ws.Range("H4").Value = 0.5280659999999999
$ws.Range("I4").Value = 0.03293066697281707
$ws.Range("J4").Value = 0.03293066697281707
$ws.Range("M4").Value = 0.3699833333333333
$ws.Range("N4").Value = 1.10995
$ws.Range("O4").Value = 0.005784234991056675
$ws.Range("P4").Value = 0.005784234991056675
$ws.Range("Q4").Value = 0.06512520629999999
$ws.Range("R4").Value = 0.5861268566999999
$ws.Range("S4").Value = 0.0001904787161830029
$ws.Range("T4").Value = 0.0001904787161830029
# Row 5
$ws.Range("I5").Value = 0.8002039325901205
$ws.Range("J5").Value = 0.8002039325901203
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 63.55492266666666
$ws.Range("N5").Value = 190.664768
$ws.Range("O5").Value = 0.9936031556622397
$ws.Range("P5").Value = 0.9936031556622397
$ws.Range("Q5").Value = 271.8418394791182
$ws.Range("R5").Value = 2446.576555312064
$ws.Range("S5").Value = 0.7950851525948778
$ws.Range("T5").Value = 0.7950851525948777
# Row 6
$ws.Range("I6").Value = 0.8002039325901205
$ws.Range("J6").Value = 0.8002039325901203
$ws.Range("O6").Value = 0.000612609346703606
$ws.Range("P6").Value = 0.000612609346703606
$ws.Range("S6").Value = 0.0004902124083736901
$ws.Range("T6").Value = 0.00049021240837369
# Row 7
$ws.Range("I7").Value = 0.8002039325901205
$ws.Range("J7").Value = 0.8002039325901203
$ws.Range("M7").Value = 0.3699833333333333
$ws.Range("N7").Value = 1.10995
$ws.Range("O7").Value = 0.005784234991056675
$ws.Range("P7").Value = 0.005784234991056675
$ws.Range("Q7").Value = 1.582520215427778
$ws.Range("R7").Value = 14.24268193885
$ws.Range("S7").Value = 0.004628567586868932
$ws.Range("T7").Value = 0.004628567586868931
# Row 8
$ws.Range("G8").Value = 0.891934
$ws.Range("H8").Value = 2.675802
$ws.Range("I8").Value = 0.1668654004370625
$ws.Range("J8").Value = 0.1668654004370625
$ws.Range("K8").Value = 3
$ws.Range("L8").Value = 1
$ws.Range("M8").Value = 63.55492266666666
$ws.Range("N8").Value = 190.664768
$ws.Range("O8").Value = 0.9936031556622397
$ws.Range("P8").Value = 0.9936031556622397
$ws.Range("Q8").Value = 56.68679639377066
$ws.Range("R8").Value = 510.1811675439359
$ws.Range("S8").Value = 0.1657979884451086
$ws.Range("T8").Value = 0.1657979884451086
# Row 9
$ws.Range("G9").Value = 0.891934
$ws.Range("H9").Value = 2.675802
$ws.Range("I9").Value = 0.1668654004370625
$ws.Range("J9").Value = 0.1668654004370625
$ws.Range("O9").Value = 0.000612609346703606
$ws.Range("P9").Value = 0.000612609346703606
$ws.Range("Q9").Value = 0.03495043379000001
$ws.Range("R9").Value = 0.31455390411
$ws.Range("S9").Value = 0.0001022233039491845
$ws.Range("T9").Value = 0.0001022233039491845
# Row 10
$ws.Range("G10").Value = 0.891934
$ws.Range("H10").Value = 2.675802
$ws.Range("I10").Value = 0.1668654004370625
$ws.Range("J10").Value = 0.1668654004370625
$ws.Range("M10").Value = 0.3699833333333333
$ws.Range("N10").Value = 1.10995
$ws.Range("O10").Value = 0.005784234991056675
$ws.Range("P10").Value = 0.005784234991056675
$ws.Range("Q10").Value = 0.3300007144333333
$ws.Range("R10").Value = 2.9700064299
$ws.Range("S10").Value = 0.000965188688004741
$ws.Range("T10").Value = 0.000965188688004741
